$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 32: dodonam.tistory.com post update
$ws.Range("D32").Value = "리스트 정렬 (multiple key를 이용한 정렬)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/337"

# Row 36: dmqm seminar update
$ws.Range("D36").Value = "Introduction to Hyper-Parameter Optimization"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/336"

# Row 41: cloudinsight log post update
$ws.Range("D41").Value = "로그 데이터의 수집과 시각화 – Part 3"
$ws.Range("E41").Value = "http://cloudinsight.net/data/log-part-3/"
